# Generate Report for handoff
# Two previously "Ready for handoff" source files (09deae5b..., 1afa181a...)
# have moved on to "In Translation", and two brand new source files
# (4b158da7-1552-4e12-a8de-56e4c78925ad, d5f46b90-cd77-4d15-9b24-2b3737a52cfd)
# have been picked up and are now "Ready for handoff".

$wb = $excel.ActiveWorkbook

$repoBase   = "https://github.com/OpenLocalizationTest/oltest/blob/e81e5378ee55d12ec63bc0ca1ae131f7d252c1de"
$zhBase     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c0146efaa9a01f7563e193326a6a985de21863e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang"
$deBase     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea0f8fe0b117791b7162d574705ad27a19946d8f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang"

$file1 = "09deae5b-c91e-4c2e-a0a0-973e2373c9cd"
$file2 = "1afa181a-bf57-46c8-9f19-7b6cce006306"
$file3 = "4b158da7-1552-4e12-a8de-56e4c78925ad"
$file4 = "d5f46b90-cd77-4d15-9b24-2b3737a52cfd"

$file1xlfZh = "$file1.b5d8a3492d5dcf5920ea9f35d311625e9f7ffd0f.zh-cn.xlf"
$file2xlfZh = "$file2.fa1ca6ac5d75f53401de665ee76de499ed9d335c.zh-cn.xlf"
$file3xlfZh = "$file3.da4919f34f69783c8a15f380d1ca65e90ec45ae0.zh-cn.xlf"
$file4xlfZh = "$file4.80d364d35e5c74f54c7dbec218a9a5fbeee319f8.zh-cn.xlf"

$file1xlfDe = "$file1.b5d8a3492d5dcf5920ea9f35d311625e9f7ffd0f.de-de.xlf"
$file2xlfDe = "$file2.fa1ca6ac5d75f53401de665ee76de499ed9d335c.de-de.xlf"
$file3xlfDe = "$file3.da4919f34f69783c8a15f380d1ca65e90ec45ae0.de-de.xlf"
$file4xlfDe = "$file4.80d364d35e5c74f54c7dbec218a9a5fbeee319f8.de-de.xlf"

$handoffDtZh34 = "2016-01-18 06:10:18"
$handoffDtDe34 = "2016-01-18 06:10:31"

function Set-DateText($range, $text) {
    $range.Value = $text
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

function Style-Link($range) {
    $range.Style = "HyperLink"
}

# ----------------------------------------------------------------------
# Sheet "Overview"
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Drop every existing hyperlink on this sheet so we can lay them back down
# cleanly in the right final order.
$ws1.Range("A1").Hyperlinks.Delete()

$ws1.Range("A1").Value = "File Name"
$ws1.Range("B1").Value = "zh-cn"
$ws1.Range("C1").Value = "de-de"

$ws1.Range("A2").Value = "$file1.md"
$ws1.Range("B2").Value = "In Translation"
$ws1.Range("C2").Value = "In Translation"

$ws1.Range("A3").Value = "$file2.md"
$ws1.Range("B3").Value = "In Translation"
$ws1.Range("C3").Value = "In Translation"

$ws1.Range("A4").Value = "$file3.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

$ws1.Range("A5").Value = "$file4.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"

$ws1.Range("A6").Value = ".localization-config"
$ws1.Range("B6").Value = "Not to be localized"
$ws1.Range("C6").Value = "Not to be localized"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "$repoBase/e2e/$file1.md", "", "", "$file1.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$repoBase/e2e/$file2.md", "", "", "$file2.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$repoBase/e2e/$file3.md", "", "", "$file3.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "$repoBase/e2e/$file4.md", "", "", "$file4.md")
$ws1.Hyperlinks.Add($ws1.Range("A6"), "$repoBase/.localization-config", "", "", ".localization-config")

Style-Link $ws1.Range("A2")
Style-Link $ws1.Range("A3")
Style-Link $ws1.Range("A4")
Style-Link $ws1.Range("A5")
Style-Link $ws1.Range("A6")

# ----------------------------------------------------------------------
# Sheet "zh-cn"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A1").Hyperlinks.Delete()

$ws2.Range("A1").Value = "Source File Name"
$ws2.Range("B1").Value = "Status"
$ws2.Range("C1").Value = "Latest Handoff File"
$ws2.Range("D1").Value = "Latest Handoff Datetime"
$ws2.Range("E1").Value = "Latest Target File"
$ws2.Range("F1").Value = "Latest Handback File"
$ws2.Range("G1").Value = "Latest Handback DateTime"
$ws2.Range("H1").Value = "Handoff Reason"
$ws2.Range("I1").Value = "Dependency From"

$ws2.Range("A2").Value = "$file1.md"
$ws2.Range("B2").Value = "In Translation"
$ws2.Range("C2").Value = $file1xlfZh
Set-DateText $ws2.Range("D2") "2016-01-18 06:08:26"
Set-DateText $ws2.Range("G2") "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = "$file2.md"
$ws2.Range("B3").Value = "In Translation"
$ws2.Range("C3").Value = $file2xlfZh
Set-DateText $ws2.Range("D3") "2016-01-18 06:08:26"
Set-DateText $ws2.Range("G3") "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = "$file3.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = $file3xlfZh
Set-DateText $ws2.Range("D4") $handoffDtZh34
Set-DateText $ws2.Range("G4") "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Include"

$ws2.Range("A5").Value = "$file4.md"
$ws2.Range("B5").Value = "Ready for handoff"
$ws2.Range("C5").Value = $file4xlfZh
Set-DateText $ws2.Range("D5") $handoffDtZh34
Set-DateText $ws2.Range("G5") "0001-01-01 00:00:00"
$ws2.Range("H5").Value = "Include"

$ws2.Range("A6").Value = ".localization-config"
$ws2.Range("B6").Value = "Not to be localized"
Set-DateText $ws2.Range("D6") "0001-01-01 00:00:00"
Set-DateText $ws2.Range("G6") "0001-01-01 00:00:00"
$ws2.Range("H6").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "$repoBase/e2e/$file1.md", "", "", "$file1.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$zhBase/$file1xlfZh", "", "", $file1xlfZh)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$repoBase/e2e/$file2.md", "", "", "$file2.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "$zhBase/$file2xlfZh", "", "", $file2xlfZh)
$ws2.Hyperlinks.Add($ws2.Range("A4"), "$repoBase/e2e/$file3.md", "", "", "$file3.md")
$ws2.Hyperlinks.Add($ws2.Range("C4"), "$zhBase/$file3xlfZh", "", "", $file3xlfZh)
$ws2.Hyperlinks.Add($ws2.Range("A5"), "$repoBase/e2e/$file4.md", "", "", "$file4.md")
$ws2.Hyperlinks.Add($ws2.Range("C5"), "$zhBase/$file4xlfZh", "", "", $file4xlfZh)
$ws2.Hyperlinks.Add($ws2.Range("A6"), "$repoBase/.localization-config", "", "", ".localization-config")

Style-Link $ws2.Range("A2")
Style-Link $ws2.Range("C2")
Style-Link $ws2.Range("A3")
Style-Link $ws2.Range("C3")
Style-Link $ws2.Range("A4")
Style-Link $ws2.Range("C4")
Style-Link $ws2.Range("A5")
Style-Link $ws2.Range("C5")
Style-Link $ws2.Range("A6")

# ----------------------------------------------------------------------
# Sheet "de-de"
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A1").Hyperlinks.Delete()

$ws3.Range("A1").Value = "Source File Name"
$ws3.Range("B1").Value = "Status"
$ws3.Range("C1").Value = "Latest Handoff File"
$ws3.Range("D1").Value = "Latest Handoff Datetime"
$ws3.Range("E1").Value = "Latest Target File"
$ws3.Range("F1").Value = "Latest Handback File"
$ws3.Range("G1").Value = "Latest Handback DateTime"
$ws3.Range("H1").Value = "Handoff Reason"
$ws3.Range("I1").Value = "Dependency From"

$ws3.Range("A2").Value = "$file1.md"
$ws3.Range("B2").Value = "In Translation"
$ws3.Range("C2").Value = $file1xlfDe
Set-DateText $ws3.Range("D2") "2016-01-18 06:08:45"
Set-DateText $ws3.Range("G2") "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = "$file2.md"
$ws3.Range("B3").Value = "In Translation"
$ws3.Range("C3").Value = $file2xlfDe
Set-DateText $ws3.Range("D3") "2016-01-18 06:08:45"
Set-DateText $ws3.Range("G3") "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = "$file3.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = $file3xlfDe
Set-DateText $ws3.Range("D4") $handoffDtDe34
Set-DateText $ws3.Range("G4") "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Include"

$ws3.Range("A5").Value = "$file4.md"
$ws3.Range("B5").Value = "Ready for handoff"
$ws3.Range("C5").Value = $file4xlfDe
Set-DateText $ws3.Range("D5") $handoffDtDe34
Set-DateText $ws3.Range("G5") "0001-01-01 00:00:00"
$ws3.Range("H5").Value = "Include"

$ws3.Range("A6").Value = ".localization-config"
$ws3.Range("B6").Value = "Not to be localized"
Set-DateText $ws3.Range("D6") "0001-01-01 00:00:00"
Set-DateText $ws3.Range("G6") "0001-01-01 00:00:00"
$ws3.Range("H6").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "$repoBase/e2e/$file1.md", "", "", "$file1.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$deBase/$file1xlfDe", "", "", $file1xlfDe)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$repoBase/e2e/$file2.md", "", "", "$file2.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "$deBase/$file2xlfDe", "", "", $file2xlfDe)
$ws3.Hyperlinks.Add($ws3.Range("A4"), "$repoBase/e2e/$file3.md", "", "", "$file3.md")
$ws3.Hyperlinks.Add($ws3.Range("C4"), "$deBase/$file3xlfDe", "", "", $file3xlfDe)
$ws3.Hyperlinks.Add($ws3.Range("A5"), "$repoBase/e2e/$file4.md", "", "", "$file4.md")
$ws3.Hyperlinks.Add($ws3.Range("C5"), "$deBase/$file4xlfDe", "", "", $file4xlfDe)
$ws3.Hyperlinks.Add($ws3.Range("A6"), "$repoBase/.localization-config", "", "", ".localization-config")

Style-Link $ws3.Range("A2")
Style-Link $ws3.Range("C2")
Style-Link $ws3.Range("A3")
Style-Link $ws3.Range("C3")
Style-Link $ws3.Range("A4")
Style-Link $ws3.Range("C4")
Style-Link $ws3.Range("A5")
Style-Link $ws3.Range("C5")
Style-Link $ws3.Range("A6")
